# Insert a new weekly price-report row for "Berenjena" (Región Metropolitana)
# just after the current last row of data (row 298), i.e. at row 299,
# pushing all the existing rows (299-328) down by one (300-329).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 299 - this shifts rows 299:328 down to 300:329
# and expands the used range / dimension to A1:R329.
$ws.Rows(299).Insert()

# Populate the newly inserted row 299 with the new record.
$ws.Cells.Item(299, 1).Value  = 6
$ws.Cells.Item(299, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(299, 3).Value  = "Metropolitana"
$ws.Cells.Item(299, 4).Value  = 45106
$ws.Cells.Item(299, 5).Value  = 13
$ws.Cells.Item(299, 6).Value  = 100112001
$ws.Cells.Item(299, 7).Value  = "Berenjena"
$ws.Cells.Item(299, 8).Value  = "Sin especificar"
$ws.Cells.Item(299, 9).Value  = "Primera"
$ws.Cells.Item(299, 10).Value = 480
$ws.Cells.Item(299, 11).Value = 4000
$ws.Cells.Item(299, 12).Value = 5800
$ws.Cells.Item(299, 13).Value = 4938
$ws.Cells.Item(299, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(299, 15).Value = "Región Metropolitana"
$ws.Cells.Item(299, 16).Value = 82
$ws.Cells.Item(299, 17).Value = 60
$ws.Cells.Item(299, 18).Value = "Hortaliza"
